$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.716.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "'2.479.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'319.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "'92.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").Value = "'0.0867"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.86%  "
$ws.Range("D11").Value = "'33.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'2.862.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'6.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "'15.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").Value = "'2.464.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "'0.795"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").Value = "'41.639.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "'6.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'0.0₃0946"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "'70.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'11.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "'240.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").Value = "'1.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'25.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.32%  "
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").Value = "'9.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").Value = "'37.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.89%  "
$ws.Range("D31").Value = "'157.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").Value = "'5.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "'0.0767"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").Value = "'17.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'0.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.92%  "
$ws.Range("D38").Value = "'1.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.52%  "
$ws.Range("D39").Value = "'2.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("D41").Value = "'4.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("D43").Value = "'1.998.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").Value = "'18.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("D46").Value = "'2.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.00%  "
$ws.Range("D47").Value = "'9.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.40%  "
$ws.Range("D48").Value = "'2.715.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'98.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("D50").Value = "'75.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.68%  "
$ws.Range("D51").Value = "'67.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.04%  "
